$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "{'C': 10, 'degree': 3, 'gamma': 'scale', 'kernel': 'poly'}"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "68.52%"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "69.58%"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "68.52%"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "68.54%"
